$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 02:14:53"
$ws1.Range("A3").Value = "Total filas: 3"

# Row 6
$ws1.Range("A6").Value = "02:14:53"
$ws1.Range("B6").Value = "02:57"
$ws1.Range("C6").Value = "215_ALUAR"
$ws1.Range("D6").Value = 43

# Row 7
$ws1.Range("A7").Value = "02:14:53"
$ws1.Range("B7").Value = "03:48"
$ws1.Range("C7").Value = "14_ABASTO"
$ws1.Range("D7").Value = 94

# Row 8
$ws1.Range("A8").Value = "02:14:53"
$ws1.Range("B8").Value = "04:01"
$ws1.Range("C8").Value = "81_EL PELIGRO"
$ws1.Range("D8").Value = 107

# Row 9 no longer exists - remove it
$ws1.Rows.Item(9).Delete()

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 02:14:53"

$ws2.Range("A6").Value = "02:14:53"
$ws2.Range("B6").Value = "02:57"
$ws2.Range("D6").Value = 43

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 02:14:53"
